# Update countries & provincias Spain
# Applies the data refresh captured by the commit:
#  - Estados Unidos (row 4), Brasil (row 14), Rumania (row 36),
#    Kenia/Jamaica (rows 118-119, which also swap display order),
#    and Ruanda (row 132) get refreshed case numbers.
#  - The "datos actualizados" timestamp moves from 00:52 to 01:22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4, 2).Value = 1059992
$ws.Cells.Item(4, 3).Value = 24227
$ws.Cells.Item(4, 4).Value = 145904
$ws.Cells.Item(4, 5).Value = 852580
$ws.Cells.Item(4, 6).Value = 18819
$ws.Cells.Item(4, 7).Value = 2242
$ws.Cells.Item(4, 8).Value = 61508

# --- Row 14: Brasil ---
$ws.Cells.Item(14, 2).Value = 79218
$ws.Cells.Item(14, 3).Value = 6319
$ws.Cells.Item(14, 5).Value = 39579
$ws.Cells.Item(14, 7).Value = 444
$ws.Cells.Item(14, 8).Value = 5507

# --- Row 36: Rumania ---
$ws.Cells.Item(36, 5).Value = 7716
$ws.Cells.Item(36, 7).Value = 30
$ws.Cells.Item(36, 8).Value = 693

# --- Rows 118-119: Jamaica now listed before Kenia, with refreshed totals ---
$ws.Cells.Item(118, 1).Value = "Jamaica"
$ws.Cells.Item(118, 2).Value = 396
$ws.Cells.Item(118, 3).Value = 32
$ws.Cells.Item(118, 4).Value = 29
$ws.Cells.Item(118, 5).Value = 360
$ws.Cells.Item(118, 6).Value = 3
$ws.Cells.Item(118, 7).Value = 0
$ws.Cells.Item(118, 8).Value = 7

$ws.Cells.Item(119, 1).Value = "Kenia"
$ws.Cells.Item(119, 2).Value = 384
$ws.Cells.Item(119, 3).Value = 10
$ws.Cells.Item(119, 4).Value = 129
$ws.Cells.Item(119, 5).Value = 240
$ws.Cells.Item(119, 6).Value = 2
$ws.Cells.Item(119, 7).Value = 1
$ws.Cells.Item(119, 8).Value = 15

# --- Row 132: Ruanda ---
$ws.Cells.Item(132, 2).Value = 225
$ws.Cells.Item(132, 3).Value = 13
$ws.Cells.Item(132, 4).Value = 98
$ws.Cells.Item(132, 5).Value = 127

# --- Timestamp cell (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 01:22"
